$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 09:20 AM"

$wsIA = $wb.Worksheets.Item("Industry Analysis")
$wsIA.Range("F2").Value = 21.0016
$wsIA.Range("F3").Value = -16.2396
$wsIA.Range("F4").Value = 27.1317
$wsIA.Range("F5").Value = -50.6494
$wsIA.Range("F6").Value = 53.2813
$wsIA.Range("F7").Value = -8.106199999999999
$wsIA.Range("F8").Value = -9.552099999999999
$wsIA.Range("F9").Value = 36.3756
$wsIA.Range("F10").Value = -6.1314
$wsIA.Range("F11").Value = 31.9081
$wsIA.Range("F12").Value = -18.4955
$wsIA.Range("F13").Value = 14.0155
$wsIA.Range("F14").Value = -36.0718
$wsIA.Range("F15").Value = -0.1622
$wsIA.Range("F16").Value = 0.1459
$wsIA.Range("F17").Value = -22.0012
$wsIA.Range("F18").Value = 1.0561
$wsIA.Range("F19").Value = -27.708
$wsIA.Range("F20").Value = 47.7309
$wsIA.Range("F21").Value = 12.0959
$wsIA.Range("F22").Value = 95.1491
$wsIA.Range("F23").Value = -50.2657
$wsIA.Range("F24").Value = -13.3427
$wsIA.Range("F25").Value = -9.9316
$wsIA.Range("F26").Value = 5.8244
$wsIA.Range("F27").Value = -32.7692
$wsIA.Range("F28").Value = -24.8224
$wsIA.Range("F29").Value = -18.4191
$wsIA.Range("F30").Value = 25.8569
$wsIA.Range("F31").Value = 58.4712
$wsIA.Range("F32").Value = -3.3862
$wsIA.Range("F33").Value = -6.3282
$wsIA.Range("F34").Value = 27.7203
$wsIA.Range("F35").Value = 4.4873
$wsIA.Range("F36").Value = -4.9458
$wsIA.Range("F37").Value = 3.6074
$wsIA.Range("F38").Value = -23.3973
$wsIA.Range("F39").Value = 8.7355
$wsIA.Range("F40").Value = -5.8541
$wsIA.Range("F41").Value = -8.3934
$wsIA.Range("F42").Value = 20.3818
$wsIA.Range("F43").Value = 14.3164
$wsIA.Range("F44").Value = -12.6846
$wsIA.Range("F45").Value = 28.4075
$wsIA.Range("F46").Value = -1.1135
$wsIA.Range("F47").Value = -37.1997
$wsIA.Range("F48").Value = -29.8569
$wsIA.Range("F49").Value = -27.5511
$wsIA.Range("F50").Value = -49.7478
$wsIA.Range("F51").Value = -51.8002
$wsIA.Range("F52").Value = -38.5254
$wsIA.Range("F53").Value = -12.4886
$wsIA.Range("F54").Value = -5.0725
$wsIA.Range("F55").Value = -17.7445
$wsIA.Range("F56").Value = -26.636
$wsIA.Range("F57").Value = -29.3361
$wsIA.Range("F58").Value = -11.9574
$wsIA.Range("F59").Value = -24.5687
$wsIA.Range("F60").Value = -12.3
$wsIA.Range("F61").Value = -10.9446
$wsIA.Range("F62").Value = -17.1229
$wsIA.Range("F63").Value = -9.5038
$wsIA.Range("F64").Value = 54.2749
$wsIA.Range("F65").Value = -43.4736
$wsIA.Range("F66").Value = 13.2687
$wsIA.Range("F67").Value = 12.7149
$wsIA.Range("F68").Value = 24.8057
$wsIA.Range("F69").Value = -17.0328
$wsIA.Range("F70").Value = -6.8927
$wsIA.Range("F71").Value = 13.6034
$wsIA.Range("F72").Value = 3.9995
$wsIA.Range("F73").Value = -16.226
$wsIA.Range("F74").Value = -16.2448
$wsIA.Range("F75").Value = 28.6924
$wsIA.Range("F76").Value = 48.9752
